# Update crypto price/volume table cells per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a "plain" number (single decimal point, no thousands
# separators) must be forced to Text format first, otherwise Excel auto-converts
# the assigned string into a numeric value and the literal formatting (e.g. trailing
# zeros like "1.00") would be lost.
$textCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D17', 'D19', 'D20', 'D23', 'D24', 'D25', 'D27', 'D28', 'D30', 'D31', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D41', 'D42', 'D44', 'D46', 'D48', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '48.392.45'
$ws.Range('E2').Value = '  +2.01%  '

$ws.Range('D3').Value = '2.516.50'
$ws.Range('E3').Value = '  +0.15%  '

$ws.Range('D4').Value = '1.00'

$ws.Range('D5').Value = '323.67'
$ws.Range('E5').Value = '  -0.29%  '

$ws.Range('D6').Value = '109.49'
$ws.Range('E6').Value = '  -0.63%  '

$ws.Range('D7').Value = '0.526'
$ws.Range('E7').Value = '  -0.26%  '

$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.13%  '

$ws.Range('D9').Value = '0.563'
$ws.Range('E9').Value = '  +3.60%  '

$ws.Range('D10').Value = '40.55'
$ws.Range('E10').Value = '  +3.15%  '

$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.0820'
$ws.Range('E11').Value = '  +0.03%  '

$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').Value = '19.66'
$ws.Range('E12').Value = '  +5.41%  '

$ws.Range('E13').Value = '  +0.63%  '

$ws.Range('D14').Value = '7.23'
$ws.Range('E14').Value = '  -0.32%  '

$ws.Range('D15').Value = '2.905.71'
$ws.Range('E15').Value = '  +0.12%  '

$ws.Range('D16').Value = '2.515.11'
$ws.Range('E16').Value = '  +0.17%  '

$ws.Range('D17').Value = '0.855'
$ws.Range('E17').Value = '  -1.18%  '

$ws.Range('D18').Value = '48.197.70'
$ws.Range('E18').Value = '  +1.66%  '

$ws.Range('D19').Value = '13.41'
$ws.Range('E19').Value = '  +2.97%  '

$ws.Range('D20').Value = '6.64'
$ws.Range('E20').Value = '  -1.39%  '

$ws.Range('D21').Value = '0.0₃0946'
$ws.Range('E21').Value = '  -0.20%  '

$ws.Range('E22').Value = '  +4.06%  '

$ws.Range('D23').Value = '71.67'
$ws.Range('E23').Value = '  +0.99%  '

$ws.Range('D24').Value = '271.12'
$ws.Range('E24').Value = '  +8.28%  '

$ws.Range('D25').Value = '2.56'
$ws.Range('E25').Value = '  -1.87%  '

$ws.Range('E26').Value = '  +0.04%  '

$ws.Range('D27').Value = '26.12'
$ws.Range('E27').Value = '  -0.61%  '

$ws.Range('D28').Value = '10.21'
$ws.Range('E28').Value = '  +1.28%  '

$ws.Range('E29').Value = '  -0.68%  '

$ws.Range('D30').Value = '0.143'
$ws.Range('E30').Value = '  +4.84%  '

$ws.Range('D31').Value = '35.37'
$ws.Range('E31').Value = '  -1.52%  '

$ws.Range('E32').Value = '  -1.02%  '

$ws.Range('D33').Value = '20.10'
$ws.Range('E33').Value = '  +0.62%  '

$ws.Range('D34').Value = '5.41'
$ws.Range('E34').Value = '  -1.13%  '

$ws.Range('D35').Value = '1.01'
$ws.Range('E35').Value = '  +0.10%  '

$ws.Range('D36').Value = '0.0789'
$ws.Range('E36').Value = '  -0.97%  '

$ws.Range('D37').Value = '1.99'
$ws.Range('E37').Value = '  -1.49%  '

$ws.Range('D38').Value = '4.71'
$ws.Range('E38').Value = '  -0.90%  '

$ws.Range('D39').Value = '2.98'
$ws.Range('E39').Value = '  -1.09%  '

$ws.Range('E40').Value = '  -0.38%  '

$ws.Range('D41').Value = '22.22'

$ws.Range('D42').Value = '119.06'
$ws.Range('E42').Value = '  -3.10%  '

$ws.Range('E43').Value = '  -3.77%  '

$ws.Range('D44').Value = '0.0300'
$ws.Range('E44').Value = '  +0.23%  '

$ws.Range('D45').Value = '2.002.01'
$ws.Range('E45').Value = '  -0.03%  '

$ws.Range('D46').Value = '3.12'
$ws.Range('E46').Value = '  +0.24%  '

$ws.Range('D48').Value = '1.85'
$ws.Range('E48').Value = '  +3.16%  '

$ws.Range('E49').Value = '  +0.18%  '

$ws.Range('D50').Value = '5.24'
$ws.Range('E50').Value = '  -1.30%  '

$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '79.88'
$ws.Range('E51').Value = '  +1.34%  '
